$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: re-"type" the full text of a paragraph (search for its current
# exact text and replace it with the same text). Word's engine rewrites the
# run(s) backing that text as a single homogeneous run when the replacement
# is performed this way, which is exactly the "multiple <w:r> merged into
# one <w:r>" normalization seen throughout the diff for paragraphs whose
# wording did not actually change.
# ---------------------------------------------------------------------------
function Normalize-Paragraph($index) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $txt = $r.Text
    # Trim the trailing paragraph mark / cell mark characters Word appends
    # to Range.Text so Find.Execute matches only the visible content.
    $txt = $txt.TrimEnd([char]13, [char]7)
    $r.Find.Execute($txt, $true, $false, $false, $false, $false, $true, 1, $false, $txt, 2) | Out-Null
}

# 1) "pick up" / "집어들다. 찾다" -> single run
Normalize-Paragraph 2
# 2) "serve" / "제공하다." -> single run
Normalize-Paragraph 3
# 3) "medicine" / "약" -> single run
Normalize-Paragraph 4
# 4) "take" / "먹다" -> single run
Normalize-Paragraph 5

# ---------------------------------------------------------------------------
# New vocabulary entry "Magnificent" inserted right after "accept", using
# the blank paragraph that already separates the vocab list from the
# example sentences. A fresh blank paragraph is added back afterwards so
# the example-sentence block keeps its leading blank line.
# ---------------------------------------------------------------------------
$blank = $d.Paragraphs.Item(15)
$br = $blank.Range
$br.InsertAfter("Magnificent" + "`t`t`t`t" + "훌륭한, 참으로 아름다운")
$br.InsertParagraphAfter() | Out-Null

# After the insertion the paragraph numbering shifts down by one for
# everything that follows.

# "That coffee shop serves good coffee" / "그 커피숌의 커피 맛이 정말 좋던데요." -> single run
Normalize-Paragraph 18
# "Yes, take it after a meal." / "네, 식사 후에 복용하세요." -> single run
Normalize-Paragraph 20

# "I'm back to + V" line: the trailing "~" run and the "하러 왔습니다." run
# get merged into one run reading "~하러 왔습니다.". Only the affected tail
# (not the whole paragraph, which contains a straight apostrophe in "I'm"
# that autocorrect would otherwise mangle into a curly quote when retyped)
# is searched/replaced.
$p22 = $d.Paragraphs.Item(22)
$r22 = $p22.Range
$r22.Find.Execute("~하러 왔습니다.", $true, $false, $false, $false, $false, $true, 1, $false, "~하러 왔습니다.", 2) | Out-Null

# " Can I pay for this with my credit card?" / "신용카드로 지불할 수 있나요?" -> single run
Normalize-Paragraph 27
